$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels with units
$ws.Range("J1").Value = 'MAE [$COP/kWh]'
$ws.Range("K1").Value = 'MSE [$COP/kWh]'
$ws.Range("L1").Value = 'RMSE [$COP/kWh]'
$ws.Range("M1").Value = 'MAPE [%]'

# Update data row values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 25
$ws.Range("G2").Value = "<keras.src.optimizers.adam.Adam object at 0x0000022183D41780>"
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 27.60541312062069
$ws.Range("K2").Value = 1222.927920048971
$ws.Range("L2").Value = 34.97038632970718
$ws.Range("M2").Value = 16.74087729137624
